$d = $word.ActiveDocument

# 1) Remove the empty "First Paragraph" (just a <w:br/>) that sits right
#    after the "23 My Atrium Patient Portal" heading and before the
#    "Critical to good communication..." bullet list.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "First Paragraph" -and $p.Range.Text.Trim().Length -eq 0) {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -like "Critical to good communication*") {
            $target = $p
            break
        }
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 2) Text replacements
$d.Content.Find.Execute(
    "Critical to good communication with your cancer care team", $true, $false, $false, $false, $false,
    $true, 1, $false, "Critical to good communication with your care team", 2)

$d.Content.Find.Execute(
    "Important to reduce the risk of complications from cancer treatment", $true, $false, $false, $false, $false,
    $true, 1, $false, "Reduces risk of complications from treatment", 2)

$d.Content.Find.Execute(
    "Working hard enough that you can" + [char]0x2019 + "t carry a conversation", $true, $false, $false, $false, $false,
    $true, 1, $false, "Working hard enough that you can" + [char]0x2019 + "t converse", 2)

$d.Content.Find.Execute(
    "Start slow an build up", $true, $false, $false, $false, $false,
    $true, 1, $false, "Start slowly and build up", 2)

$d.Content.Find.Execute(
    "Smoking makes it more difficult to get through cancer treatment", $true, $false, $false, $false, $false,
    $true, 1, $false, "Smoking makes cancer treatment more difficult", 2)

$d.Content.Find.Execute(
    "American Lung Asssociation fredomfromsmoking.org", $true, $false, $false, $false, $false,
    $true, 1, $false, "American Lung Assn fredomfromsmoking.org", 2)

$d.Content.Find.Execute(
    "1:1 Smoking Cessation Counseling Clinics (Metro Charlotte)", $true, $false, $false, $false, $false,
    $true, 1, $false, "1:1 Smoking Cessation Counseling (Metro Charlotte)", 2)
